# Generate Report for Handback
#
# This script mirrors a localization "handback" run: the target markdown
# file has come back from translation for zh-cn and de-de, so each
# language sheet gets its "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns populated (instead of the empty /
# zero-date placeholders) and the overall Status flips from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

$targetFile   = "ad4d7d2a-fd6b-4e1a-8731-a99ffd00814c.md"
$targetUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c119663a9a25334ced9b47a795f429d63fffe2bc/e2e/ad4d7d2a-fd6b-4e1a-8731-a99ffd00814c.md"
$zhcnHandback = "ad4d7d2a-fd6b-4e1a-8731-a99ffd00814c.bc5f2a3527db211a72d80c53ff71854b0207bb3f.zh-cn.xlf"
$dedeHandback = "ad4d7d2a-fd6b-4e1a-8731-a99ffd00814c.bc5f2a3527db211a72d80c53ff71854b0207bb3f.de-de.xlf"

$zhcnHandbackDate = "2016-09-04 23:07:55"
$dedeHandbackDate = "2016-09-04 23:08:07"

$hyperlinkColor = 15570276   # RGB(100,149,237) == FF6495ED, matches the workbook's existing HyperLink style

# ---------------------------------------------------------------------
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview mirrors both languages in columns E/F, each language sheet
#    carries its own Status in column C)
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (I), Latest Handback File (J)
#    and Latest Handback DateTime (K) for both data rows.
# ---------------------------------------------------------------------
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $targetUrl, "", "", $targetFile) | Out-Null
$zhcn.Range("I2").Font.Underline = $true
$zhcn.Range("I2").Font.Color = $hyperlinkColor

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $targetUrl, "", "", $targetFile) | Out-Null
$zhcn.Range("I3").Font.Underline = $true
$zhcn.Range("I3").Font.Color = $hyperlinkColor

$zhcn.Range("J2").Value = $zhcnHandback
$zhcn.Range("J3").Value = $zhcnHandback

$zhcn.Range("K2").Value = $zhcnHandbackDate
$zhcn.Range("K3").Value = $zhcnHandbackDate

# ---------------------------------------------------------------------
# 3. de-de sheet: same treatment, with de-de specific handback file name
#    and its own handback datetime.
# ---------------------------------------------------------------------
$dede.Hyperlinks.Add($dede.Range("I2"), $targetUrl, "", "", $targetFile) | Out-Null
$dede.Range("I2").Font.Underline = $true
$dede.Range("I2").Font.Color = $hyperlinkColor

$dede.Hyperlinks.Add($dede.Range("I3"), $targetUrl, "", "", $targetFile) | Out-Null
$dede.Range("I3").Font.Underline = $true
$dede.Range("I3").Font.Color = $hyperlinkColor

$dede.Range("J2").Value = $dedeHandback
$dede.Range("J3").Value = $dedeHandback

$dede.Range("K2").Value = $dedeHandbackDate
$dede.Range("K3").Value = $dedeHandbackDate

# ---------------------------------------------------------------------
# 4. Widen the columns that now hold the longer Status text / populated
#    file names & links so the report reads cleanly.
# ---------------------------------------------------------------------
$wideStatus = 29.9777047293527 - (5/6)   # -> stored col width ~30
$wideFile   = 40 - (5/6)                 # -> stored col width 40

$overview.Columns.Item(5).ColumnWidth = $wideStatus   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = $wideStatus   # column F (de-de status)

$zhcn.Columns.Item(3).ColumnWidth = $wideStatus        # column C (Status)
$zhcn.Columns.Item(9).ColumnWidth = $wideFile          # column I (Latest Target File)
$zhcn.Columns.Item(10).ColumnWidth = $wideFile         # column J (Latest Handback File)

$dede.Columns.Item(3).ColumnWidth = $wideStatus        # column C (Status)
$dede.Columns.Item(9).ColumnWidth = $wideFile          # column I (Latest Target File)
$dede.Columns.Item(10).ColumnWidth = $wideFile         # column J (Latest Handback File)

Write-Output "Handback report generated"
